# Buildings workbook update:
#  - Replace the single "Threshold" header (G1) with four new headers:
#      G1 = Threshold (Summer), H1 = Threshold (Winter),
#      I1 = Threshold (Spring), J1 = Battery Capacity (Spring & Spoof)
#  - Populate the new G:J columns with per-building numeric data (rows 2-11)
#  - Add SUM totals in row 12 for the new columns (bold, thousands-format,
#    matching the existing B12:F12 total cells)
#  - Resize the new threshold columns (G:I) to fit their header text
#  - Leave final selection on F13 (matches the author's saved cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row -----------------------------------------------------------
# Note: shared-string table order follows first-write order, so write
# Winter/Summer/Spring/Battery in the same sequence the source workbook used
# (H1 before G1) to keep the sharedStrings table laid out identically.
$ws.Range("H1").Value = "Threshold (Winter)"
$ws.Range("H1").Font.Bold = $true

$ws.Range("G1").Value = "Threshold (Summer)"

$ws.Range("I1").Value = "Threshold (Spring)"
$ws.Range("I1").Font.Bold = $true

$ws.Range("J1").Value = "Battery Capacity (Spring & Spoof)"
$ws.Range("J1").Font.Bold = $true

# ---- Per-building data (rows 2-11) -----------------------------------------
$data = @{
    2  = @(60059, 60208, 55208, 403393)
    3  = @(60059, 60208, 55208, 403393)
    4  = @(60059, 60208, 55208, 403393)
    5  = @(36035, 36124, 33365, 242036)
    6  = @(36035, 36124, 33365, 242036)
    7  = @(36035, 36124, 33365, 242036)
    8  = @(84083, 84291, 77852, 564741)
    9  = @(88087, 88305, 81559, 591643)
    10 = @(74740, 74925, 69201, 502000)
    11 = @(101433, 101685, 93916, 681286)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 7).Value = $vals[0]
    $ws.Cells.Item($row, 8).Value = $vals[1]
    $ws.Cells.Item($row, 9).Value = $vals[2]
    $ws.Cells.Item($row, 10).Value = $vals[3]
}

# ---- Totals row (row 12) ---------------------------------------------------
$ws.Range("G12:J12").FormulaR1C1 = "=SUM(R[-10]C:R[-1]C)"
$ws.Range("G12:J12").NumberFormat = "#,##0"
$ws.Range("G12:J12").Font.Bold = $true

# ---- Column widths for the new Threshold columns ---------------------------
# (18 characters matches the length of the "Threshold (....)" header text;
# offset compensates for this engine's fixed column-width padding constant
# so the serialized <col width="..."/> lands on 18, same as native Excel's
# bestFit calculation for these headers.)
$ws.Range("G1:I1").ColumnWidth = 17.1666666666667

# ---- Final selection --------------------------------------------------------
$ws.Range("F13").Select()
